$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 304
$ws.Range("I12").Value2 = 305.66666
$ws.Range("J12").Value2 = 299
$ws.Range("K12").Value2 = 305.66666
$ws.Range("L12").Value2 = 299
$ws.Range("M12").Value2 = -135.66666
$ws.Range("N12").Value2 = -639
$ws.Range("H43").Value2 = 3468.6667
$ws.Range("J43").Value2 = 3623
$ws.Range("L43").Value2 = 3623
$ws.Range("N43").Value2 = -3761
$ws.Range("H138").Value2 = 3032.3489
$ws.Range("I138").Value2 = 1812.4546
$ws.Range("K138").Value2 = 5437.3638
$ws.Range("M138").Value2 = -297.3638000000001
$ws.Range("H139").Value2 = 200000
$ws.Range("J139").Value2 = 200000
$ws.Range("L139").Value2 = 200000
$ws.Range("N139").Value2 = -210280
$ws.Range("H141").Value2 = 3969.4285
$ws.Range("I141").Value2 = 3760.9473
$ws.Range("J141").Value2 = 5950
$ws.Range("K141").Value2 = 11282.8419
$ws.Range("L141").Value2 = 17850
$ws.Range("M141").Value2 = -6102.841899999999
$ws.Range("N141").Value2 = -28210

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2509.5898
$ws.Range("I32").Value2 = 2638.2
$ws.Range("K32").Value2 = 2638.2
$ws.Range("M32").Value2 = -2351.2
$ws.Range("H97").Value2 = 753.9048
$ws.Range("I97").Value2 = 522.4706
$ws.Range("K97").Value2 = 522.4706
$ws.Range("M97").Value2 = -26.47059999999999
$ws.Range("H122").Value2 = 2830.8
$ws.Range("I122").Value2 = 2618.8
$ws.Range("J122").Value2 = 3042.8
$ws.Range("K122").Value2 = 7856.400000000001
$ws.Range("L122").Value2 = 9128.400000000001
$ws.Range("M122").Value2 = -5406.400000000001
$ws.Range("N122").Value2 = -14028.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value2 = 6370.2
$ws.Range("I94").Value2 = 5802.6665
$ws.Range("K94").Value2 = 5802.6665
$ws.Range("M94").Value2 = -5351.6665

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 356.7143
$ws.Range("I22").Value2 = 349.8
$ws.Range("K22").Value2 = 349.8
$ws.Range("M22").Value2 = 0.1999999999999886
$ws.Range("H31").Value2 = 1926.3636
$ws.Range("I31").Value2 = 1739.4
$ws.Range("K31").Value2 = 1739.4
$ws.Range("M31").Value2 = -1444.4
$ws.Range("H34").Value2 = 1926.3636
$ws.Range("I34").Value2 = 1739.4
$ws.Range("K34").Value2 = 1739.4
$ws.Range("M34").Value2 = -1537.4
$ws.Range("H58").Value2 = 3040.4614
$ws.Range("I58").Value2 = 2929.6365
$ws.Range("K58").Value2 = 2929.6365
$ws.Range("M58").Value2 = -2726.6365
$ws.Range("H94").Value2 = 1428.5
$ws.Range("I94").Value2 = 1500
$ws.Range("J94").Value2 = 1418.2858
$ws.Range("K94").Value2 = 1500
$ws.Range("L94").Value2 = 1418.2858
$ws.Range("N94").Value2 = -2320.2858
$ws.Range("M94").Value2 = -1049
$ws.Range("H99").Value2 = 2338.2104
$ws.Range("I99").Value2 = 1759.5834
$ws.Range("J99").Value2 = 3330.1428
$ws.Range("K99").Value2 = 1759.5834
$ws.Range("L99").Value2 = 3330.1428
$ws.Range("M99").Value2 = -261.5834
$ws.Range("N99").Value2 = -6326.1428
$ws.Range("H126").Value2 = 2338.2104
$ws.Range("I126").Value2 = 1759.5834
$ws.Range("J126").Value2 = 3330.1428
$ws.Range("K126").Value2 = 5278.7502
$ws.Range("L126").Value2 = 9990.428400000001
$ws.Range("M126").Value2 = -2808.7502
$ws.Range("N126").Value2 = -14930.4284
$ws.Range("H134").Value2 = 7146230
$ws.Range("I134").Value2 = 3438.182
$ws.Range("J134").Value2 = 33336466
$ws.Range("K134").Value2 = 10314.546
$ws.Range("L134").Value2 = 100009398
$ws.Range("M134").Value2 = -7779.545999999998
$ws.Range("N134").Value2 = -100014468
$ws.Range("H136").Value2 = 3040.4614
$ws.Range("I136").Value2 = 2929.6365
$ws.Range("K136").Value2 = 8788.9095
$ws.Range("M136").Value2 = -6238.9095
$ws.Range("H141").Value2 = 99990.836
$ws.Range("I141").Value2 = 0
$ws.Range("J141").Value2 = 99990.836
$ws.Range("K141").Value2 = 0
$ws.Range("L141").Value2 = 99990.836
$ws.Range("M141").ClearContents() | Out-Null
$ws.Range("N141").Value2 = -110350.836

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1075
$ws.Range("I5").Value2 = 1033.8125
$ws.Range("K5").Value2 = 3101.4375
$ws.Range("M5").Value2 = -2989.4375
$ws.Range("H9").Value2 = 125002330
$ws.Range("I9").Value2 = 1217
$ws.Range("J9").Value2 = 200002990
$ws.Range("K9").Value2 = 3651
$ws.Range("L9").Value2 = 600008970
$ws.Range("M9").Value2 = -3427
$ws.Range("N9").Value2 = -600009418
$ws.Range("H17").Value2 = 346.33334
$ws.Range("I17").Value2 = 150
$ws.Range("J17").Value2 = 444.5
$ws.Range("K17").Value2 = 450
$ws.Range("L17").Value2 = 1333.5
$ws.Range("M17").Value2 = -281
$ws.Range("N17").Value2 = -1671.5
$ws.Range("H132").Value2 = 3542.1667
$ws.Range("J132").Value2 = 3542.1667
$ws.Range("L132").Value2 = 31879.5003
$ws.Range("N132").Value2 = -36939.5003
$ws.Range("H133").Value2 = 9139.333000000001
$ws.Range("I133").Value2 = 9139.333000000001
$ws.Range("K133").Value2 = 27417.999
$ws.Range("M133").Value2 = -22357.999
$ws.Range("H135").Value2 = 1075
$ws.Range("I135").Value2 = 1033.8125
$ws.Range("K135").Value2 = 9304.3125
$ws.Range("M135").Value2 = -6769.3125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value2 = 26666
$ws.Range("H70").Value2 = 11408.909
$ws.Range("I70").Value2 = 10833.111
$ws.Range("J70").Value2 = 14000
$ws.Range("K70").Value2 = 10833.111
$ws.Range("L70").Value2 = 14000
$ws.Range("M70").Value2 = -10563.111
$ws.Range("N70").Value2 = -14540
$ws.Range("H73").Value2 = 11408.909
$ws.Range("I73").Value2 = 10833.111
$ws.Range("J73").Value2 = 14000
$ws.Range("K73").Value2 = 10833.111
$ws.Range("L73").Value2 = 14000
$ws.Range("M73").Value2 = -9897.111000000001
$ws.Range("N73").Value2 = -15872
$ws.Range("H80").Value2 = 8699.799999999999
$ws.Range("I80").Value2 = 4250
$ws.Range("K80").Value2 = 4250
$ws.Range("M80").Value2 = -3252
$ws.Range("H83").Value2 = 8699.799999999999
$ws.Range("I83").Value2 = 4250
$ws.Range("K83").Value2 = 21250
$ws.Range("M83").Value2 = -16258
$ws.Range("H97").Value2 = 519.5714
$ws.Range("I97").Value2 = 669.1667
$ws.Range("J97").Value2 = 320.1111
$ws.Range("K97").Value2 = 669.1667
$ws.Range("L97").Value2 = 320.1111
$ws.Range("M97").Value2 = -173.1667
$ws.Range("N97").Value2 = -1312.1111
$ws.Range("H113").Value2 = 508.5
$ws.Range("I113").Value2 = 496.66666
$ws.Range("K113").Value2 = 496.66666
$ws.Range("M113").Value2 = 1673.33334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 2107.5
$ws.Range("I22").Value2 = 1099.4
$ws.Range("K22").Value2 = 1099.4
$ws.Range("M22").Value2 = -804.4000000000001
$ws.Range("H27").Value2 = 2107.5
$ws.Range("I27").Value2 = 1099.4
$ws.Range("K27").Value2 = 1099.4
$ws.Range("M27").Value2 = -992.4000000000001
$ws.Range("H40").Value2 = 4585
$ws.Range("I40").Value2 = 4199.5
$ws.Range("J40").Value2 = 4739.2
$ws.Range("K40").Value2 = 4199.5
$ws.Range("L40").Value2 = 4739.2
$ws.Range("M40").Value2 = -4063.5
$ws.Range("N40").Value2 = -5011.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 4646.8
$ws.Range("I62").Value2 = 4933.5
$ws.Range("K62").Value2 = 4933.5
$ws.Range("M62").Value2 = -4309.5
$ws.Range("H65").Value2 = 4646.8
$ws.Range("I65").Value2 = 4933.5
$ws.Range("K65").Value2 = 24667.5
$ws.Range("M65").Value2 = -21547.5
$ws.Range("H96").Value2 = 3369.3125
$ws.Range("I96").Value2 = 5826.125
$ws.Range("J96").Value2 = 912.5
$ws.Range("K96").Value2 = 5826.125
$ws.Range("L96").Value2 = 912.5
$ws.Range("M96").Value2 = -4453.125
$ws.Range("N96").Value2 = -3658.5
$ws.Range("H100").Value2 = 1179.1428
$ws.Range("I100").Value2 = 1131.875
$ws.Range("J100").Value2 = 1242.1666
$ws.Range("K100").Value2 = 2263.75
$ws.Range("L100").Value2 = 2484.3332
$ws.Range("M100").Value2 = -1722.75
$ws.Range("N100").Value2 = -3566.3332
$ws.Range("H107").Value2 = 1142.9286
$ws.Range("I107").Value2 = 1154
$ws.Range("J107").Value2 = 999
$ws.Range("K107").Value2 = 3462
$ws.Range("L107").Value2 = 2997
$ws.Range("M107").Value2 = -1542
$ws.Range("N107").Value2 = -6837
$ws.Range("H122").Value2 = 6198.4
$ws.Range("I122").Value2 = 6249
$ws.Range("J122").Value2 = 5996
$ws.Range("K122").Value2 = 18747
$ws.Range("L122").Value2 = 17988
$ws.Range("M122").Value2 = -16297
$ws.Range("N122").Value2 = -22888
$ws.Range("H136").Value2 = 1452.7188
$ws.Range("I136").Value2 = 1346.0769
$ws.Range("J136").Value2 = 1914.8334
$ws.Range("K136").Value2 = 4038.2307
$ws.Range("L136").Value2 = 5744.5002
$ws.Range("M136").Value2 = -1488.2307
$ws.Range("N136").Value2 = -10844.5002
